$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet and turn off right-to-left display
$ws.Name = "Sheet1"
$excel.ActiveWindow.DisplayRightToLeft = $false

# Populate feedback row
$ws.Range("A2").Value = "Chatbot"
$ws.Range("B2").Value = "It was good expirience"
$ws.Range("C2").Value = "Chatbot"

# Select A1 as the active cell
$ws.Range("A1").Select()
